$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

# Row 2
$ws.Cells.Item(2, 4).Value = '46.043.63'
$ws.Cells.Item(2, 5).Value = '  -0.93%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '2.598.59'
$ws.Cells.Item(3, 5).Value = '  -0.41%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  +0.05%  '

# Row 5
Set-TextCell 5 4 '310.90'
$ws.Cells.Item(5, 5).Value = '  +0.90%  '

# Row 6
Set-TextCell 6 4 '98.59'
$ws.Cells.Item(6, 5).Value = '  -2.74%  '

# Row 7
$ws.Cells.Item(7, 5).Value = '  -0.87%  '

# Row 8
$ws.Cells.Item(8, 5).Value = '  +0.05%  '

# Row 9
Set-TextCell 9 4 '0.581'
$ws.Cells.Item(9, 5).Value = '  +0.58%  '

# Row 10
Set-TextCell 10 4 '38.96'
$ws.Cells.Item(10, 5).Value = '  +0.13%  '

# Row 11
Set-TextCell 11 4 '54.44'
$ws.Cells.Item(11, 5).Value = '  -1.78%  '

# Row 12
Set-TextCell 12 4 '0.0839'
$ws.Cells.Item(12, 5).Value = '  -0.20%  '

# Row 13
Set-TextCell 13 4 '8.13'
$ws.Cells.Item(13, 5).Value = '  -1.22%  '

# Row 14
$ws.Cells.Item(14, 4).Value = '2.999.06'
$ws.Cells.Item(14, 5).Value = '  -0.35%  '

# Row 15
$ws.Cells.Item(15, 5).Value = '  +1.12%  '

# Row 16
$ws.Cells.Item(16, 4).Value = '2.595.76'
$ws.Cells.Item(16, 5).Value = '  -0.49%  '

# Row 17
Set-TextCell 17 4 '0.916'
$ws.Cells.Item(17, 5).Value = '  +1.20%  '

# Row 18
Set-TextCell 18 4 '14.85'
$ws.Cells.Item(18, 5).Value = '  -0.24%  '

# Row 19
$ws.Cells.Item(19, 4).Value = '46.202.71'
$ws.Cells.Item(19, 5).Value = '  -0.91%  '

# Row 20
$ws.Cells.Item(20, 5).Value = '  +0.67%  '

# Row 21
Set-TextCell 21 4 '12.80'
$ws.Cells.Item(21, 5).Value = '  -4.28%  '

# Row 22
Set-TextCell 22 4 '6.70'
$ws.Cells.Item(22, 5).Value = '  -0.06%  '

# Row 23
Set-TextCell 23 4 '295.86'
$ws.Cells.Item(23, 5).Value = '  +14.61%  '

# Row 24
Set-TextCell 24 4 '72.80'
$ws.Cells.Item(24, 5).Value = '  +2.11%  '

# Row 25
Set-TextCell 25 4 '3.06'
$ws.Cells.Item(25, 5).Value = '  +1.31%  '

# Row 26
Set-TextCell 26 4 '2.25'
$ws.Cells.Item(26, 5).Value = '  +0.73%  '

# Row 27
Set-TextCell 27 4 '29.74'
$ws.Cells.Item(27, 5).Value = '  +4.58%  '

# Row 28
$ws.Cells.Item(28, 5).Value = '  +0.22%  '

# Row 29
$ws.Cells.Item(29, 5).Value = '  +1.06%  '

# Row 30
Set-TextCell 30 4 '10.77'

# Row 31
Set-TextCell 31 4 '38.32'
$ws.Cells.Item(31, 5).Value = '  -4.54%  '

# Row 32
$ws.Cells.Item(32, 5).Value = '  -2.72%  '

# Row 33
Set-TextCell 33 4 '6.25'
$ws.Cells.Item(33, 5).Value = '  +0.94%  '

# Row 34
Set-TextCell 34 4 '3.57'
$ws.Cells.Item(34, 5).Value = '  -4.64%  '

# Row 35
Set-TextCell 35 4 '155.41'

# Row 36
Set-TextCell 36 4 '0.0838'
$ws.Cells.Item(36, 5).Value = '  +0.37%  '

# Row 37
Set-TextCell 37 4 '2.20'
$ws.Cells.Item(37, 5).Value = '  -5.36%  '

# Row 38
Set-TextCell 38 4 '2.78'
$ws.Cells.Item(38, 5).Value = '  -5.79%  '

# Row 39
$ws.Cells.Item(39, 5).Value = '  +3.49%  '

# Row 40
$ws.Cells.Item(40, 5).Value = '  +0.66%  '

# Row 41
$ws.Cells.Item(41, 2).Value = 'Celestia'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
Set-TextCell 41 4 '15.80'
$ws.Cells.Item(41, 5).Value = '  -0.52%  '

# Row 42
$ws.Cells.Item(42, 2).Value = 'VeChain'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextCell 42 4 '0.0331'
$ws.Cells.Item(42, 5).Value = '  +2.16%  '

# Row 43
$ws.Cells.Item(43, 2).Value = 'NEARProtocol'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextCell 43 4 '3.59'
$ws.Cells.Item(43, 5).Value = '  -1.67%  '

# Row 44
$ws.Cells.Item(44, 2).Value = 'EnergySwap'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextCell 44 4 '21.36'
$ws.Cells.Item(44, 5).Value = '  +12.89%  '

# Row 45
Set-TextCell 45 4 '3.94'
$ws.Cells.Item(45, 5).Value = '  -5.89%  '

# Row 46
$ws.Cells.Item(46, 4).Value = '2.098.07'
$ws.Cells.Item(46, 5).Value = '  +2.25%  '

# Row 47
Set-TextCell 47 4 '97.91'
$ws.Cells.Item(47, 5).Value = '  +7.21%  '

# Row 48
Set-TextCell 48 4 '0.999'
$ws.Cells.Item(48, 5).Value = '  +0.03%  '

# Row 49
$ws.Cells.Item(49, 5).Value = '  +3.82%  '

# Row 50
$ws.Cells.Item(50, 5).Value = '  +0.26%  '

# Row 51
Set-TextCell 51 4 '108.16'
$ws.Cells.Item(51, 5).Value = '  -1.64%  '
